# Commit: "Simplify a bit further"
#
# 1. Fill in the speaker notes on slide 1 (title slide) with the
#    introductory talk-track text (previously an empty paragraph).
# 2. Remove the "Outline" slide (originally slide 2) from the deck.
#    Every slide after it shifts up by one position as a natural
#    consequence of the deletion.

$p = $ppt.ActivePresentation

# --- 1. Speaker notes for slide 1 -----------------------------------
$slide1 = $p.Slides.Item(1)
$notes1 = $slide1.NotesPage
$notesBody = $notes1.Shapes.Item(2)

$apostrophe = [char]0x2019
$notesText = "I" + $apostrophe + "ll start with some motivation for fast and programmable networks and then talk about my past, current, and future work in this space."
$notesBody.TextFrame.TextRange.Text = $notesText

# --- 2. Delete the "Outline" slide (slide 2) -------------------------
$p.Slides.Item(2).Delete()
